$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates: rename "Unidades/Hora" and add new column F ---
$ws.Range("E1").Value = "Tiempo Alistamiento"
$ws.Range("F1").Value = "Tiempo Operación"

# --- Turn off gridlines for the window/view ---
$excel.ActiveWindow.DisplayGridlines = $false

# --- Center-align the header row A1:F1 (keeps existing bold/fill header style) ---
$ws.Range("A1:F1").HorizontalAlignment = -4108

# --- Column widths (approximate Excel "characters" widths) ---
$ws.Columns.Item(1).ColumnWidth = 14
$ws.Columns.Item(2).ColumnWidth = 33
$ws.Columns.Item(3).ColumnWidth = 20.5
$ws.Columns.Item(4).ColumnWidth = 17.83333333333333
$ws.Columns.Item(5).ColumnWidth = 23.16666666666667
$ws.Columns.Item(6).ColumnWidth = 20.5
